$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from the existing last header cell (AC1) onto
# the three new header cells so they pick up the same cell style (bold,
# centered, bordered) as the rest of row 1 instead of an ad-hoc new style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-50) gets the same team record: 82 wins, 79 losses, 0 ties.
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 82  # AD
    $ws.Cells.Item($row, 31).Value = 79  # AE
    $ws.Cells.Item($row, 32).Value = 0   # AF
}
